# DC-Colos data update: insert a new "TEN" (Tongren, China) colo row just
# before the existing "IAD" row (old row 278), shifting every subsequent
# row down by one. Dimension grows from A1:H338 to A1:H339.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 278 (pushes IAD.. down to 279..339)
$ws.Rows.Item(278).Insert()

# The freshly inserted row has no explicit style yet; copy the "colo"
# column's formatting (bold / centered / thin-bordered) from the row right
# below it (the shifted-down IAD row) so the new A278 cell matches the
# rest of column A.
$ws.Cells.Item(279, 1).Copy()
$ws.Cells.Item(278, 1).PasteSpecial(-4122)

# Populate the new row's values.
$ws.Cells.Item(278, 1).Value = "TEN"
$ws.Cells.Item(278, 2).Value = "Tongren, China"
$ws.Cells.Item(278, 3).Value = "Asia"
$ws.Cells.Item(278, 4).Value = "Tongren"
$ws.Cells.Item(278, 5).Value = "China"
$ws.Cells.Item(278, 6).Value = "CN"
# lat/lon are left blank for this colo (no coordinates available yet).
